# Apply bug-fix edits to Raven9Q_test_2.xlsx
# - question_answers sheet: column B (per-question answers) values shuffled/corrected
# - outputs sheet: recomputed raw score, percentile ranks, and IQ scores

$wb = $excel.ActiveWorkbook

$wsAnswers = $wb.Worksheets.Item("question_answers")
$wsOutputs = $wb.Worksheets.Item("outputs")

# --- question_answers: column B corrected values (row -> new value) ---
$answerUpdates = @{
    2  = "4"
    4  = "1"
    5  = "2"
    6  = "6"
    7  = "5"
    8  = "4"
    9  = "5"
    10 = "6"
    11 = "1"
    12 = "3"
    13 = "1"
    14 = "3"
    15 = "5"
    16 = "1"
    17 = "4"
    18 = "4"
    19 = "3"
    20 = "4"
    21 = "1"
    23 = "1"
    24 = "3"
    25 = "1"
    26 = "3"
    27 = "4"
    28 = "4"
    29 = "1"
    30 = "6"
    31 = "2"
    33 = "3"
    34 = "3"
    35 = "2"
    36 = "2"
    37 = "2"
}

foreach ($row in $answerUpdates.Keys) {
    $cell = $wsAnswers.Range("B$row")
    $cell.NumberFormat = "@"
    $cell.Value = $answerUpdates[$row]
}

# --- outputs: recomputed scores ---
# B2 (raw) stays numeric, as in the original workbook.
$wsOutputs.Range("B2").Value = 27

# B4/B5/B6 (percentile ranks) and B8/B9 (iq scores) are stored as text.
$textCells = @("B4", "B5", "B6", "B8", "B9")
foreach ($addr in $textCells) {
    $cell = $wsOutputs.Range($addr)
    $cell.NumberFormat = "@"
}

$wsOutputs.Range("B4").Value = "95+"
$wsOutputs.Range("B5").Value = "95+"
$wsOutputs.Range("B6").Value = "95+"
$wsOutputs.Range("B8").Value = "130"
$wsOutputs.Range("B9").Value = "135"
